# Update "想去人数" (want-to-go count) figures in column F across all
# sheets, refreshing the generated report snapshot.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 209
$ws.Range("F4").Value = 13236
$ws.Range("F5").Value = 54498
$ws.Range("F7").Value = 1321
$ws.Range("F8").Value = 352
$ws.Range("F9").Value = 314
$ws.Range("F10").Value = 869
$ws.Range("F11").Value = 732
$ws.Range("F12").Value = 375
$ws.Range("F13").Value = 3006
$ws.Range("F14").Value = 883
$ws.Range("F15").Value = 5191
$ws.Range("F16").Value = 1268
$ws.Range("F17").Value = 973
$ws.Range("F19").Value = 560
$ws.Range("F21").Value = 388
$ws.Range("F22").Value = 1239
$ws.Range("F23").Value = 89
$ws.Range("F24").Value = 36
$ws.Range("F25").Value = 166
$ws.Range("F26").Value = 346
$ws.Range("F27").Value = 9
$ws.Range("F31").Value = 53
$ws.Range("F32").Value = 4849
$ws.Range("F34").Value = 4839
$ws.Range("F35").Value = 8813
$ws.Range("F36").Value = 111
$ws.Range("F39").Value = 210
$ws.Range("F40").Value = 419
$ws.Range("F41").Value = 106
$ws.Range("F43").Value = 4181
$ws.Range("F44").Value = 214

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4176
$ws.Range("F4").Value = 93
$ws.Range("F5").Value = 130
$ws.Range("F7").Value = 54
$ws.Range("F20").Value = 89

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 782
$ws.Range("F3").Value = 560

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 782
$ws.Range("F3").Value = 560
$ws.Range("F4").Value = 209
$ws.Range("F6").Value = 13236
$ws.Range("F7").Value = 1321
$ws.Range("F8").Value = 352
$ws.Range("F9").Value = 314
$ws.Range("F10").Value = 869
$ws.Range("F11").Value = 732
$ws.Range("F12").Value = 375
$ws.Range("F13").Value = 3006
$ws.Range("F14").Value = 883
$ws.Range("F15").Value = 93
$ws.Range("F16").Value = 1268
$ws.Range("F18").Value = 130
$ws.Range("F19").Value = 973
$ws.Range("F20").Value = 54
$ws.Range("F21").Value = 560
$ws.Range("F22").Value = 388
$ws.Range("F24").Value = 1239
$ws.Range("F26").Value = 166
$ws.Range("F28").Value = 346
$ws.Range("F30").Value = 53
$ws.Range("F31").Value = 4849
$ws.Range("F33").Value = 4839
$ws.Range("F34").Value = 8813
$ws.Range("F35").Value = 111
$ws.Range("F38").Value = 210
$ws.Range("F39").Value = 419
$ws.Range("F42").Value = 106
$ws.Range("F44").Value = 4181
$ws.Range("F47").Value = 214
